$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated sample data rows (rows 4-9), leaving only the
# header rows (1-2) and the template row (3). Deleting these rows also
# drops the now-unused shared strings that were only referenced there.
$ws.Rows("4:9").Delete()
